$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the "open site url" step: row 14 becomes "click contact us",
#     row 15 becomes "enter corect data", row 16 becomes "click send message" /
#     "message about successfully sent message", and the old row 17
#     (duplicate "click send message" / success message) is deleted outright.

# Re-point the field values first (top of the form).
$ws.Range("B5").Value = "Successfully sent message via contact form"
$ws.Range("B8").Value = "opened homepage"
$ws.Range("B9").Value = "mladenowa_tedi@abv.bg"

# The Input Data cell (B9) no longer links out to the test site -- drop the
# hyperlink so the plain e-mail text is left behind.
$ws.Hyperlinks.Delete()

# Row 14 used to be the italic/gray placeholder step "1. open site url" with
# its expected result; it now carries the old row-15 "click contact us" step
# text and no expected-result note, styled the same as the rows below it.
$ws.Range("A15").Copy()
$ws.Range("A14:B14").PasteSpecial(-4122)
$ws.Range("A14").Value = "click contact us"
$ws.Range("B14").Value = ""

# Row 15 shifts to "enter corect data" with no expected-result note.
$ws.Range("A15").Value = "enter corect data"
$ws.Range("B15").Value = ""

# Row 16 shifts to "click send message" / the success message that used to
# live on row 17.
$ws.Range("A16").Value = "click send message"
$ws.Range("B16").Value = "message about successfully sent message"

# The old row 17 (formerly "click send message" / success message) is now
# redundant -- remove it so the table ends at row 16.
$ws.Rows(17).Delete()

$ws.Range("B16").Select()
